# Remove the "Ver no Jupiter..." line, the "(c) 2020 ..." footer line, and
# the blank paragraph that separates them from the preceding
# "LOB1008: ..." requisito paragraph. These three consecutive paragraphs
# are deleted in their entirety (including their paragraph marks), leaving
# the trailing blank paragraph (and the page-break paragraph after it)
# untouched.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $i
    }
    if ($text -like "*Contact: luizeleno@usp.br*") {
        $endPara = $i
    }
}

# Also swallow the blank paragraph immediately before the "Ver no Jupiter"
# paragraph so it disappears along with the other two.
$deleteStart = $startPara
if ($deleteStart -gt 1) {
    $prevText = $d.Paragraphs.Item($deleteStart - 1).Range.Text
    if ($prevText -eq [string][char]13) {
        $deleteStart = $deleteStart - 1
    }
}

$range = $d.Range($d.Paragraphs.Item($deleteStart).Range.Start, $d.Paragraphs.Item($endPara).Range.End)
$range.Delete()
